$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Scratch cell (far outside the real A1:E37 data range) used to push string
# (text) values through the clipboard so Excel stores them as genuine
# shared-string text in the destination cell -- a plain ".Value = '2018'"
# assignment gets auto-coerced to a number by Excel, which is not what we
# want since cohort_year must stay text (matches the original column type).
# The helper cell holds a "=TEXT(...)" formula (a string-valued formula
# result, not a literal), so no NumberFormat/style change is needed on it;
# PasteSpecial(values) onto a General-formatted destination then keeps the
# destination cell's style at the default (no explicit style index) while
# still writing it as real shared-string text, exactly like the other
# untouched data cells.
$helperCol = 30
$helper = $ws.Cells.Item(1, $helperCol)

function Set-TextValue($targetCell, [string]$text) {
    $f = '=TEXT("' + $text + '","@")'
    $helper.Formula = $f
    $helper.Copy()
    $targetCell.PasteSpecial(-4163)
}

# Full refreshed cohort-retention table (row, cohort_year, period_index,
# num_customers, cohort_size, retention_rate).
$data = @(
    ,@(2, "2018", 0, 1006, 1006, 1)
    ,@(3, "2018", 1, 109, 1006, 0.1083499005964215)
    ,@(4, "2018", 2, 38, 1006, 0.03777335984095427)
    ,@(5, "2018", 3, 5, 1006, 0.004970178926441352)
    ,@(6, "2018", 4, 7, 1006, 0.006958250497017893)
    ,@(7, "2018", 5, 7, 1006, 0.006958250497017893)
    ,@(8, "2018", 6, 4, 1006, 0.003976143141153081)
    ,@(9, "2018", 7, 2, 1006, 0.001988071570576541)
    ,@(10, "2019", 0, 2127, 2127, 1)
    ,@(11, "2019", 1, 132, 2127, 0.06205923836389281)
    ,@(12, "2019", 2, 37, 2127, 0.01739539257169723)
    ,@(13, "2019", 3, 34, 2127, 0.01598495533615421)
    ,@(14, "2019", 4, 28, 2127, 0.01316408086506817)
    ,@(15, "2019", 5, 22, 2127, 0.01034320639398213)
    ,@(16, "2019", 6, 6, 2127, 0.002820874471086037)
    ,@(17, "2020", 0, 2659, 2659, 1)
    ,@(18, "2020", 1, 173, 2659, 0.06506205340353516)
    ,@(19, "2020", 2, 129, 2659, 0.04851447912749154)
    ,@(20, "2020", 3, 77, 2659, 0.02895825498307635)
    ,@(21, "2020", 4, 58, 2659, 0.02181271154569387)
    ,@(22, "2020", 5, 21, 2659, 0.007897705904475368)
    ,@(23, "2021", 0, 2278, 2278, 1)
    ,@(24, "2021", 1, 274, 2278, 0.1202809482001756)
    ,@(25, "2021", 2, 129, 2278, 0.05662862159789289)
    ,@(26, "2021", 3, 94, 2278, 0.04126426690079017)
    ,@(27, "2021", 4, 28, 2278, 0.01229148375768218)
    ,@(28, "2022", 0, 2317, 2317, 1)
    ,@(29, "2022", 1, 243, 2317, 0.1048769961156668)
    ,@(30, "2022", 2, 133, 2317, 0.05740181268882175)
    ,@(31, "2022", 3, 30, 2317, 0.01294777729823047)
    ,@(32, "2023", 0, 2256, 2256, 1)
    ,@(33, "2023", 1, 202, 2256, 0.08953900709219859)
    ,@(34, "2023", 2, 48, 2256, 0.02127659574468085)
    ,@(35, "2024", 0, 1932, 1932, 1)
    ,@(36, "2024", 1, 90, 1932, 0.04658385093167702)
    ,@(37, "2025", 0, 463, 463, 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $yearText = $row[1]
    $cellA = $ws.Cells.Item($r, 1)
    Set-TextValue $cellA $yearText
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# drop the scratch column entirely so it leaves no trace in the saved file
$helper.EntireColumn.Delete()
